$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2528603
$ws.Range("J17").Value = 2607596.8
$ws.Range("L17").Value = 7822790.399999999
$ws.Range("N17").Value = -7823126.399999999

$ws.Range("H64").Value = 3068.7896
$ws.Range("J64").Value = 3109.9092
$ws.Range("L64").Value = 3109.9092
$ws.Range("N64").Value = -3605.9092

$ws.Range("H67").Value = 3068.7896
$ws.Range("J67").Value = 3109.9092
$ws.Range("L67").Value = 3109.9092
$ws.Range("N67").Value = -4825.9092

$ws.Range("H74").Value = 3204.8057
$ws.Range("I74").Value = 3894.6
$ws.Range("J74").Value = 3093.5483
$ws.Range("K74").Value = 3894.6
$ws.Range("L74").Value = 3093.5483
$ws.Range("M74").Value = -2958.6
$ws.Range("N74").Value = -4965.5483

$ws.Range("H77").Value = 3204.8057
$ws.Range("I77").Value = 3894.6
$ws.Range("J77").Value = 3093.5483
$ws.Range("K77").Value = 19473
$ws.Range("L77").Value = 15467.7415
$ws.Range("M77").Value = -14793
$ws.Range("N77").Value = -24827.7415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1717.862
$ws.Range("I61").Value = 1148
$ws.Range("J61").Value = 3902.3333
$ws.Range("K61").Value = 1148
$ws.Range("L61").Value = 3902.3333
$ws.Range("M61").Value = -936
$ws.Range("N61").Value = -4326.3333

$ws.Range("H74").Value = 8931637
$ws.Range("I74").Value = 14709035
$ws.Range("J74").Value = 2930.9092
$ws.Range("K74").Value = 14709035
$ws.Range("L74").Value = 2930.9092
$ws.Range("M74").Value = -14708161
$ws.Range("N74").Value = -4678.9092

$ws.Range("H77").Value = 8931637
$ws.Range("I77").Value = 14709035
$ws.Range("J77").Value = 2930.9092
$ws.Range("K77").Value = 73545175
$ws.Range("L77").Value = 14654.546
$ws.Range("M77").Value = -73540807
$ws.Range("N77").Value = -23390.546

$ws.Range("H102").Value = 1554.1666
$ws.Range("I102").Value = 1554.1666
$ws.Range("K102").Value = 1554.1666
$ws.Range("M102").Value = 67.83339999999998

$ws.Range("H122").Value = 2575.88
$ws.Range("I122").Value = 2099.6667
$ws.Range("J122").Value = 3290.2
$ws.Range("K122").Value = 6299.000100000001
$ws.Range("L122").Value = 9870.599999999999
$ws.Range("M122").Value = -3849.000100000001
$ws.Range("N122").Value = -14770.6

$ws.Range("H132").Value = 1781.3469
$ws.Range("I132").Value = 1325.4359
$ws.Range("J132").Value = 3559.4
$ws.Range("K132").Value = 3976.3077
$ws.Range("L132").Value = 10678.2
$ws.Range("M132").Value = -1446.3077
$ws.Range("N132").Value = -15738.2

$ws.Range("H136").Value = 1717.862
$ws.Range("I136").Value = 1148
$ws.Range("J136").Value = 3902.3333
$ws.Range("K136").Value = 3444
$ws.Range("L136").Value = 11706.9999
$ws.Range("M136").Value = -894
$ws.Range("N136").Value = -16806.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 501568.7
$ws.Range("I134").Value = 770251.9
$ws.Range("J134").Value = 2585.7144
$ws.Range("K134").Value = 2310755.7
$ws.Range("L134").Value = 7757.1432
$ws.Range("M134").Value = -2308220.7
$ws.Range("N134").Value = -12827.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8002015
$ws.Range("I31").Value = 14286716
$ws.Range("J31").Value = 3304.2727
$ws.Range("K31").Value = 14286716
$ws.Range("L31").Value = 3304.2727
$ws.Range("M31").Value = -14286421
$ws.Range("N31").Value = -3894.2727

$ws.Range("H34").Value = 8002015
$ws.Range("I34").Value = 14286716
$ws.Range("J34").Value = 3304.2727
$ws.Range("K34").Value = 14286716
$ws.Range("L34").Value = 3304.2727
$ws.Range("M34").Value = -14286514
$ws.Range("N34").Value = -3708.2727

$ws.Range("H99").Value = 1883732
$ws.Range("I99").Value = 2555007.8
$ws.Range("J99").Value = 4160
$ws.Range("K99").Value = 2555007.8
$ws.Range("L99").Value = 4160
$ws.Range("M99").Value = -2553509.8
$ws.Range("N99").Value = -7156

$ws.Range("H122").Value = 9524895
$ws.Range("J122").Value = 1324.75
$ws.Range("L122").Value = 3974.25
$ws.Range("N122").Value = -8874.25

$ws.Range("H126").Value = 1883732
$ws.Range("I126").Value = 2555007.8
$ws.Range("J126").Value = 4160
$ws.Range("K126").Value = 7665023.399999999
$ws.Range("L126").Value = 12480
$ws.Range("M126").Value = -7662553.399999999
$ws.Range("N126").Value = -17420

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1256.0385
$ws.Range("I102").Value = 1210.95
$ws.Range("J102").Value = 1406.3334
$ws.Range("K102").Value = 1210.95
$ws.Range("L102").Value = 1406.3334
$ws.Range("M102").Value = 411.05
$ws.Range("N102").Value = -4650.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1514.4166
$ws.Range("I7").Value = 1093.4615
$ws.Range("J7").Value = 1752.3478
$ws.Range("K7").Value = 1093.4615
$ws.Range("L7").Value = 1752.3478
$ws.Range("M7").Value = -981.4614999999999
$ws.Range("N7").Value = -1976.3478

$ws.Range("H68").Value = 28197194
$ws.Range("I68").Value = 56389908
$ws.Range("J68").Value = 4477.6665
$ws.Range("K68").Value = 56389908
$ws.Range("L68").Value = 4477.6665
$ws.Range("M68").Value = -56389159
$ws.Range("N68").Value = -5975.6665

$ws.Range("H71").Value = 28197194
$ws.Range("I71").Value = 56389908
$ws.Range("J71").Value = 4477.6665
$ws.Range("K71").Value = 281949540
$ws.Range("L71").Value = 22388.3325
$ws.Range("M71").Value = -281945796
$ws.Range("N71").Value = -29876.3325

$ws.Range("H82").Value = 1455.3889
$ws.Range("I82").Value = 1063.625
$ws.Range("J82").Value = 1768.8
$ws.Range("K82").Value = 1063.625
$ws.Range("L82").Value = 1768.8
$ws.Range("M82").Value = -702.625
$ws.Range("N82").Value = -2490.8

$ws.Range("H85").Value = 1455.3889
$ws.Range("I85").Value = 1063.625
$ws.Range("J85").Value = 1768.8
$ws.Range("K85").Value = 1063.625
$ws.Range("L85").Value = 1768.8
$ws.Range("M85").Value = 184.375
$ws.Range("N85").Value = -4264.8

$ws.Range("H100").Value = 1822.7273
$ws.Range("I100").Value = 1512.5
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1512.5
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -971.5
$ws.Range("N100").Value = -3082

$ws.Range("H126").Value = 1514.4166
$ws.Range("I126").Value = 1093.4615
$ws.Range("J126").Value = 1752.3478
$ws.Range("K126").Value = 3280.3845
$ws.Range("L126").Value = 5257.0434
$ws.Range("M126").Value = -810.3844999999997
$ws.Range("N126").Value = -10197.0434

$ws.Range("H132").Value = 7697479
$ws.Range("I132").Value = 14293066
$ws.Range("J132").Value = 2627.9333
$ws.Range("K132").Value = 42879198
$ws.Range("L132").Value = 7883.7999
$ws.Range("M132").Value = -42876668
$ws.Range("N132").Value = -12943.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1594.0625
$ws.Range("I122").Value = 1276.625
$ws.Range("J122").Value = 1911.5
$ws.Range("K122").Value = 3829.875
$ws.Range("L122").Value = 5734.5
$ws.Range("M122").Value = -1379.875
$ws.Range("N122").Value = -10634.5

$ws.Range("H126").Value = 3582.111
$ws.Range("I126").Value = 4259.857
$ws.Range("J126").Value = 1210
$ws.Range("K126").Value = 12779.571
$ws.Range("L126").Value = 3630
$ws.Range("M126").Value = -10309.571
$ws.Range("N126").Value = -8570
